$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a brand-new row at position 987 (shifts the existing rows 987-1072
# down to 988-1073, matching dimension A1:R1072 -> A1:R1073) and populate it
# with a new weekly price record.
$ws.Rows("987:987").Insert()

$ws.Range("A987").Value = 8
$ws.Range("B987").Value = "Terminal La Palmera de La Serena"
$ws.Range("C987").Value = "Coquimbo"
$ws.Range("D987").Value = 45013
$ws.Range("E987").Value = 4
$ws.Range("F987").Value = 100112004
$ws.Range("G987").Value = "Cebolla"
$ws.Range("H987").Value = "Sin especificar"
$ws.Range("I987").Value = "Primera"
$ws.Range("J987").Value = 2000
$ws.Range("K987").Value = 7500
$ws.Range("L987").Value = 8000
$ws.Range("M987").Value = 7750
$ws.Range("N987").Value = "$/malla 18 kilos"
$ws.Range("O987").Value = "Perú"
$ws.Range("P987").Value = 431
$ws.Range("Q987").Value = 18
$ws.Range("R987").Value = "Hortaliza"
